$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row blocks to delete (1-based worksheet row numbers), in the ORIGINAL
# (pre-edit) row numbering. Deleting from the bottom up keeps the row
# numbers for blocks not yet processed stable.
$blocksToDelete = @(
    @(211, 220),
    @(206, 206),
    @(190, 190),
    @(156, 156),
    @(102, 107)
)

foreach ($block in $blocksToDelete) {
    $startRow = $block[0]
    $endRow = $block[1]
    $range = $ws.Range("A$startRow`:A$endRow").EntireRow
    $range.Delete()
}

# Renumber column B sequentially (1..200) for the remaining data rows (2..201)
$lastRow = 201
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $r - 1
}

# Restore the view: scroll back to the top and focus the final selection
$ws.Range("A1").Select()
$ws.Range("B206:D211").Select()
$ws.Range("H12").Select()
